$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.092.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "'1.875.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'313.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.5041"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").Value = "'0.3843"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("D9").Value = "'0.08549"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.63%  "
$ws.Range("D10").Value = "'1.115"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").Value = "'41.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").Value = "'6.287"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.95%  "
$ws.Range("D13").Value = "'20.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").Value = "'1.870.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").Value = "'7.207"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("E17").Value = "  -2.50%  "
$ws.Range("D18").Value = "'91.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.83%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "'18.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "'6.092"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("D23").Value = "'28.126.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").Value = "'11.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").Value = "'2.271"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.596"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "'2.089.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'156.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'126.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.1058"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.062"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.612"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.593"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "'9.641"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02458"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06578"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2183"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.215"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.243"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.84%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6392"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.80%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'11.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("B43").Value = "InternetComputer(DFINITY)"
$ws.Range("C43").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D43").Value = "'4.903"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'13.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6013"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.282"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.664"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.992"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "'1.221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.51%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'121.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'80.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.74%  "
